# Move the last data row (524, A_SERIES_DIFFERENCE) up to row 471,
# shifting rows 471-523 down by one row (472-524).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 471
$lastRow = 524

# Capture the values of the last row (the one that will move to the top of the block).
$movedA = $ws.Cells.Item($lastRow, 1).Value2
$movedB = $ws.Cells.Item($lastRow, 2).Value2
$movedC = $ws.Cells.Item($lastRow, 3).Value2
$movedD = $ws.Cells.Item($lastRow, 4).Value2

# Shift rows firstRow..(lastRow-1) down by one row, starting from the bottom
# so we don't overwrite values before reading them.
for ($r = $lastRow - 1; $r -ge $firstRow; $r--) {
    $ws.Cells.Item($r + 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 3).Value2 = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r + 1, 4).Value2 = $ws.Cells.Item($r, 4).Value2
}

# Place the captured last-row values into the first row of the block.
$ws.Cells.Item($firstRow, 1).Value2 = $movedA
$ws.Cells.Item($firstRow, 2).Value2 = $movedB
$ws.Cells.Item($firstRow, 3).Value2 = $movedC
$ws.Cells.Item($firstRow, 4).Value2 = $movedD
